$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A106").Value = "Kappa1"
$ws.Range("B106").Value = "Test Kappa with L=1"
$ws.Range("C106").Value = "Kappa_test1"

$ws.Range("A107").Value = "Kappa2"
$ws.Range("B107").Value = "Test Kappa with L=2"
$ws.Range("C107").Value = "Kappa_test2"

$ws.Range("A108").Value = "SortinoRatio1"
$ws.Range("B108").Value = "Test SortinoRatio with group=FULL"
$ws.Range("C108").Value = "SortinoRatio_test1"

$ws.Range("A109").Value = "SortinoRatio2"
$ws.Range("B109").Value = "Test SortinoRatio with group=SUBSET"
$ws.Range("C109").Value = "SortinoRatio_test2"

$ws.Range("A110").Value = "KellyRatio1"
$ws.Range("B110").Value = "Test KellyRatio with method=half"
$ws.Range("C110").Value = "KellyRatio_test1"

$ws.Range("A111").Value = "KellyRatio2"
$ws.Range("B111").Value = "Test KellyRatio with method=full"
$ws.Range("C111").Value = "KellyRatio_test2"

$ws.Range("E100").Select()
